# Replace the text-as-shared-string percentage values in the three
# "Scores by ..." sheets with the underlying full-precision numeric values.
# (Commit message: "Update table svg's for number formats that changed")

$wb = $excel.ActiveWorkbook

# --- Scores by Spending ---
$ws = $wb.Worksheets.Item("Scores by Spending")

$ws.Range("B2").Value = 83.45539900855027
$ws.Range("C2").Value = 83.93381405396646
$ws.Range("D2").Value = 93.46009572653237
$ws.Range("E2").Value = 96.61087677671375
$ws.Range("F2").Value = 90.36945874402643

$ws.Range("B3").Value = 81.8998257021498
$ws.Range("C3").Value = 83.15528577020937
$ws.Range("D3").Value = 87.13353760737169
$ws.Range("E3").Value = 92.71820457965273
$ws.Range("F3").Value = 81.41859632428398

$ws.Range("B4").Value = 78.51885454725715
$ws.Range("C4").Value = 81.62447331528534
$ws.Range("D4").Value = 73.48420890357487
$ws.Range("E4").Value = 84.39179284814433
$ws.Range("F4").Value = 62.85765555194492

$ws.Range("B5").Value = 76.99720981240274
$ws.Range("C5").Value = 81.02784255713441
$ws.Range("D5").Value = 66.16481311032456
$ws.Range("E5").Value = 81.13395072128019
$ws.Range("F5").Value = 53.5268548869691

# --- Scores by Size ---
$ws = $wb.Worksheets.Item("Scores by Size")

$ws.Range("B2").Value = 83.82159776422071
$ws.Range("C2").Value = 83.92984341754834
$ws.Range("D2").Value = 93.55022469776569
$ws.Range("E2").Value = 96.09943667320715
$ws.Range("F2").Value = 89.88385340844357

$ws.Range("B3").Value = 83.37468376981991
$ws.Range("C3").Value = 83.86443831725629
$ws.Range("D3").Value = 93.59969459404036
$ws.Range("E3").Value = 96.7906800028675
$ws.Range("F3").Value = 90.62153518649967

$ws.Range("B4").Value = 77.746416511437
$ws.Range("C4").Value = 81.34449272598371
$ws.Range("D4").Value = 69.96336073939453
$ws.Range("E4").Value = 82.7666344526415
$ws.Range("F4").Value = 58.28600304906789

# --- Scores by Type ---
$ws = $wb.Worksheets.Item("Scores by Type")

$ws.Range("B2").Value = 83.47385187384614
$ws.Range("C2").Value = 83.89642074222549
$ws.Range("D2").Value = 93.62083003509466
$ws.Range("E2").Value = 96.58648927302872
$ws.Range("F2").Value = 90.43224369343227

$ws.Range("B3").Value = 76.95673306832398
$ws.Range("C3").Value = 80.96663632734915
$ws.Range("D3").Value = 66.54845257144746
$ws.Range("E3").Value = 80.79906211395057
$ws.Range("F3").Value = 53.67220822778149
